$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 values: A2 was wrongly set to item.id's old slot (createDate),
# restore A2 to ${item.id}; add the jxls "fn" date-format function to the
# createDate / updateDate cells (G2 / H2). B2:F2 stay as they were.
$ws.Range("A2").Value = '${item.id}'
$ws.Range("G2").Value = "`${fn.format(item.createDate, 'yyyy-MM-dd HH:mm:ss')}"
$ws.Range("H2").Value = "`${fn.format(item.updateDate, 'yyyy-MM-dd HH:mm:ss')}"

# --- Selection moved from H2 to L21 (reflects the author's last cursor
# position when the template was saved).
[void]$ws.Range("L21").Select()
